# "changed from_tree to take names of lambda. lambda now a series"
#
# The three category labels used as both column headers (row 1) and row
# labels (column A) are shortened from the "*_parameters" form to a bare
# name, since they now come straight from a pandas Series' index (the
# lambda names) instead of being explicitly suffixed strings.
#
#   ecological_parameters  -> ecology
#   mechanical_parameters  -> mechanical
#   process_parameters     -> process

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (B1:D1) and column A labels (A2:A4) share the same
# underlying strings, so update every cell that references them.
$ws.Range("B1").Value = "ecology"
$ws.Range("A2").Value = "ecology"

$ws.Range("C1").Value = "mechanical"
$ws.Range("A3").Value = "mechanical"

$ws.Range("D1").Value = "process"
$ws.Range("A4").Value = "process"

# Move the active selection to A5 (was E13).
$ws.Range("A5").Select()
